$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the touched cells keep their plain-text representation (avoid Excel
# auto-converting numeric-looking / percent-looking strings into numbers).
$targetCells = @("D2","E2","D3","E3","D4","E4","D5","E5","D6","E6","D7","E7","D8","E8","E10","D11","E11","D12","E12","D13","E13","D14","E14","D15","E15","D16","E16","D17","E17","D18","E18","D19","E19","D20","E20","D21","E21","D22","E22","D23","E23","D24","E24","D25","E25","D26","E26","D27","E27","D39","E39","D40","E40","D41","E41","D42","E42","D43","E43","D44","E44","D45","E45","D46","E46","D47","E47","D48","E48","D49","E49","D50","E50","D51","E51")
foreach ($addr in $targetCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated Price (D) / Volume(1h) (E) values from the latest scrape.
$ws.Range("D2").Value = "329.70"
$ws.Range("E2").Value = "7.05%"
$ws.Range("D3").Value = "40.58"
$ws.Range("E3").Value = "12.81%"
$ws.Range("D4").Value = "5.964"
$ws.Range("E4").Value = "16.51%"
$ws.Range("D5").Value = "0.08162"
$ws.Range("E5").Value = "6.13%"
$ws.Range("D6").Value = "4.563"
$ws.Range("E6").Value = "3.95%"
$ws.Range("D7").Value = "8.760"
$ws.Range("E7").Value = "5.66%"
$ws.Range("D8").Value = "2.011"
$ws.Range("E8").Value = "9.23%"
$ws.Range("E10").Value = "2.50%"
$ws.Range("D11").Value = "0.1307"
$ws.Range("E11").Value = "17.38%"
$ws.Range("D12").Value = "0.2002"
$ws.Range("E12").Value = "8.11%"
$ws.Range("D13").Value = "0.09352"
$ws.Range("E13").Value = "6.95%"
$ws.Range("D14").Value = "0.03442"
$ws.Range("E14").Value = "2.89%"
$ws.Range("D15").Value = "0.09641"
$ws.Range("E15").Value = "1.31%"
$ws.Range("D16").Value = "0.001323"
$ws.Range("E16").Value = "-4.40%"
$ws.Range("D17").Value = "0.006083"
$ws.Range("E17").Value = "-1.32%"
$ws.Range("D18").Value = "3.375"
$ws.Range("E18").Value = "0.31%"
$ws.Range("D19").Value = "0.3500"
$ws.Range("E19").Value = "1.61%"
$ws.Range("D20").Value = "7.715"
$ws.Range("E20").Value = "21.81%"
$ws.Range("D21").Value = "0.1448"
$ws.Range("E21").Value = "12.14%"
$ws.Range("D22").Value = "0.2448"
$ws.Range("E22").Value = "5.87%"
$ws.Range("D23").Value = "0.04432"
$ws.Range("E23").Value = "2.18%"
$ws.Range("D24").Value = "0.001253"
$ws.Range("E24").Value = "4.10%"
$ws.Range("D25").Value = "0.004393"
$ws.Range("E25").Value = "3.43%"
$ws.Range("D26").Value = "0.0001188"
$ws.Range("E26").Value = "-10.80%"
$ws.Range("D27").Value = "0.0003988"
$ws.Range("E27").Value = "37.34%"
$ws.Range("D39").Value = "0.02489"
$ws.Range("E39").Value = "19.55%"
$ws.Range("D40").Value = "0.05303"
$ws.Range("E40").Value = "8.29%"
$ws.Range("D41").Value = "0.007593"
$ws.Range("E41").Value = "0.79%"
$ws.Range("D42").Value = "0.1433"
$ws.Range("E42").Value = "6.46%"
$ws.Range("D43").Value = "0.008977"
$ws.Range("E43").Value = "4.70%"
$ws.Range("D44").Value = "0.002057"
$ws.Range("E44").Value = "-0.78%"
$ws.Range("D45").Value = "0.01055"
$ws.Range("E45").Value = "25.75%"
$ws.Range("D46").Value = "0.00006836"
$ws.Range("E46").Value = "7.88%"
$ws.Range("D47").Value = "0.00000000750"
$ws.Range("E47").Value = "-0.17%"
$ws.Range("D48").Value = "0.002897"
$ws.Range("E48").Value = "-12.21%"
$ws.Range("D49").Value = "0.001799"
$ws.Range("E49").Value = "24.53%"
$ws.Range("D50").Value = "0.00002099"
$ws.Range("E50").Value = "-0.17%"
$ws.Range("D51").Value = "0.0001999"
$ws.Range("E51").Value = "-0.17%"
